$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting (style) from the existing header cell E1,
# then set its text to "time_taken" (mirrors B1:E1 header styling, which uses
# cellXfs style index 1 - bold font, thin border, centered/top aligned).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# New data column values (plain/default formatting, same as columns B:E body cells)
$ws.Range("F2").Value = "2021-10-05 13:38:43.257282"
$ws.Range("F3").Value = "2021-10-05 13:38:43.257294"
$ws.Range("F4").Value = "2021-10-05 13:38:43.257298"
$ws.Range("F5").Value = "2021-10-05 13:38:43.257301"
$ws.Range("F6").Value = "2021-10-05 13:38:43.257304"
